$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.628.16"
$ws.Range("E2").Value = "  +1.43%  "
$ws.Range("D3").Value = "2.507.21"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'591.08"
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("D6").Value = "'174.29"
$ws.Range("E6").Value = "  -0.68%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -0.26%  "
$ws.Range("D9").Value = "2.506.29"
$ws.Range("E9").Value = "  +0.49%  "
$ws.Range("E10").Value = "  +5.13%  "
$ws.Range("E11").Value = "  -1.04%  "
$ws.Range("E12").Value = "  +1.40%  "
$ws.Range("E13").Value = "  -1.64%  "
$ws.Range("D14").Value = "2.949.52"
$ws.Range("E14").Value = "  +0.58%  "
$ws.Range("D15").Value = "'25.70"
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("D16").Value = "68.666.91"
$ws.Range("E16").Value = "  +2.00%  "
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("D18").Value = "2.504.21"
$ws.Range("E18").Value = "  -3.25%  "
$ws.Range("D19").Value = "'362.63"
$ws.Range("E19").Value = "  +3.05%  "
$ws.Range("D20").Value = "'7.53"
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").Value = "'10.89"
$ws.Range("E21").Value = "  -1.79%  "
$ws.Range("E22").Value = "  -1.91%  "
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").Value = "'70.12"
$ws.Range("E24").Value = "  -0.58%  "
$ws.Range("D25").Value = "'4.15"
$ws.Range("E25").Value = "  -2.83%  "
$ws.Range("D26").Value = "'8.90"
$ws.Range("E26").Value = "  -3.32%  "
$ws.Range("E27").Value = "  -6.48%  "
$ws.Range("D28").Value = "2.640.59"
$ws.Range("E28").Value = "  +0.56%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.38%  "
$ws.Range("D30").Value = "'511.36"
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("D31").Value = "0.0₃0876"
$ws.Range("E31").Value = "  -4.20%  "
$ws.Range("D32").Value = "'7.71"
$ws.Range("E32").Value = "  -1.52%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "'1.77"
$ws.Range("E33").Value = "  -0.71%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "'1.21"
$ws.Range("E34").Value = "  -3.82%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("D36").Value = "'162.33"
$ws.Range("E37").Value = "  -4.89%  "
$ws.Range("D38").Value = "'18.52"
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("E39").Value = "  -0.19%  "
$ws.Range("B40").Value = "USDe"
$ws.Range("C40").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("B41").Value = "ImmutableX"
$ws.Range("C41").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D41").Value = "'1.31"
$ws.Range("E41").Value = "  -2.22%  "
$ws.Range("E42").Value = "  -2.94%  "
$ws.Range("D43").Value = "'4.73"
$ws.Range("E43").Value = "  -3.04%  "
$ws.Range("E44").Value = "  -4.15%  "
$ws.Range("D45").Value = "'2.31"
$ws.Range("E45").Value = "  -4.32%  "
$ws.Range("D46").Value = "'149.95"
$ws.Range("E46").Value = "  +2.53%  "
$ws.Range("D47").Value = "'3.54"
$ws.Range("E47").Value = "  +0.59%  "
$ws.Range("D48").Value = "'0.512"
$ws.Range("E48").Value = "  -1.00%  "
$ws.Range("D49").Value = "'0.0736"
$ws.Range("E49").Value = "  -1.38%  "
$ws.Range("E50").Value = "  -4.25%  "
$ws.Range("E51").Value = "  -2.01%  "
